$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 18 ("ROUND 7 (2nd March)") forecasts, added for 1st March update.
$ws.Range("B18").Value = 10.21856
$ws.Range("C18").Value = 0.04915049
$ws.Range("D18").Value = 10.2479126841792
$ws.Range("E18").Value = 0.0300047010642495
$ws.Range("F18").Value = 10.20628
$ws.Range("G18").Value = 0.06357296
$ws.Range("H18").Value = 10.224
$ws.Range("I18").Value = 0.0397
$ws.Range("J18").Value = 10.24725
$ws.Range("K18").Value = 0.04621685
$ws.Range("L18").Value = 10.20459
$ws.Range("M18").Value = 0.04813566
$ws.Range("N18").Value = 10.185255
$ws.Range("O18").Value = 0.05916629
$ws.Range("P18").Value = 10.21345
$ws.Range("Q18").Value = 0.04888

# Match the author's final selection on the sheet.
$ws.Range("L20").Select()
